# Update column F (dSF) values on Sheet1 to reflect repulled/pushed data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    2  = -5
    4  = -2
    9  = -1
    13 = -6
    14 = -5
    17 = -3
    19 = -1
    20 = -2
    21 = -1
    24 = -9
    25 = 8
    29 = 1
    31 = 3
    32 = 2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
